$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the exposure-period time for row 3 (Albert Park, 11/2/2021)
$ws.Range("C3").Value = "9:00am - 10:15am  11/2/2021"

# Delete the two "Point Cook" rows (old rows 23 and 24); rows below shift up by two
$ws.Range("A23:D24").EntireRow.Delete()

# Rename "Melbourne" location to "Melbourne Airport" for the airport-related rows (now rows 18-20)
$ws.Range("A18").Value = "Melbourne Airport"
$ws.Range("B18").Value = "901 Frankston to Melbourne Airport bus route:  Melbourne Airport to Broadmeadows Railway Station"
$ws.Range("A19").Value = "Melbourne Airport"
$ws.Range("A20").Value = "Melbourne Airport"
